$d = $word.ActiveDocument

# Paragraphs that make up the "done" content of Level 1, Level 2 and
# Level 3 (everything except the intro paragraphs and the bold
# "Level N" headings themselves) get marked with a green highlight,
# matching the author's "All done on Level 3" commit.
$targets = @(7,8,9,10,11,12,13,14,15,17,18,19,20,22,23,24,25,26,27,28,29,30)

foreach ($i in $targets) {
    $p = $d.Paragraphs.Item($i)
    # Setting Font.HighlightColorIndex (rather than Range.HighlightColorIndex)
    # applies the highlight to both the paragraph mark run properties
    # (w:pPr/w:rPr) and every run in the paragraph (w:r/w:rPr), matching
    # wdGreen (4) == <w:highlight w:val="green"/>.
    $p.Range.Font.HighlightColorIndex = 4
}
